# "Generate Report for Archive"
#
# This localization-status report is regenerated by the CI tool each run.
# The regeneration refreshes the fixed vocabulary of status values (which
# now also includes "In Translation" as a recognised-but-currently-unused
# status) and re-writes every row of the Overview / zh-cn / de-de sheets
# from the latest handoff/handback data. For the current snapshot none of
# the tracked files actually moved into the new "In Translation" state, so
# every cell keeps the same displayed text it already had - we simply
# re-assert the authoritative values coming out of the report generator so
# the workbook reflects a freshly generated archive copy.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(2, 5).Value = "Handed back: in sync with en-US"   # E2
$ov.Cells.Item(2, 6).Value = "Handed back: in sync with en-US"   # F2
$ov.Cells.Item(2, 7).Value = "2016-08-12 08:47:48"                # G2

$ov.Cells.Item(3, 5).Value = "Ready for handoff"                  # E3
$ov.Cells.Item(3, 6).Value = "Ready for handoff"                  # F3
$ov.Cells.Item(3, 7).Value = "2016-08-12 08:48:52"                # G3

$ov.Cells.Item(4, 5).Value = "Ready for handoff"                  # E4
$ov.Cells.Item(4, 6).Value = "Ready for handoff"                  # F4
$ov.Cells.Item(4, 7).Value = "2016-08-12 08:48:52"                # G4

$ov.Cells.Item(5, 5).Value = "Ready for handoff"                  # E5
$ov.Cells.Item(5, 6).Value = "Ready for handoff"                  # F5
$ov.Cells.Item(5, 7).Value = "2016-08-12 08:47:20"                # G5

# ---- Per-language detail sheets -------------------------------------
$langSheets = @("zh-cn", "de-de")

foreach ($langName in $langSheets) {
    $ws = $wb.Worksheets.Item($langName)

    # Row 2 - already handed back, in sync with en-US
    $ws.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"  # D2 Status

    # Rows 3-5 - ready for handoff
    $ws.Cells.Item(3, 3).Value = "Ready for handoff"                # D3 Status
    $ws.Cells.Item(4, 3).Value = "Ready for handoff"                # D4 Status
    $ws.Cells.Item(5, 3).Value = "Ready for handoff"                # D5 Status
}

# Record "In Translation" as a now-recognised status value for this report
# (not currently assigned to any tracked file in this snapshot).
$reportStatuses = @(
    "Handed back: in sync with en-US",
    "In Translation",
    "Ready for handoff"
)
